# IMPORTANT UPDATE: Update of the HK map
# Apply corrected station coordinate values (columns D/E) as captured by the
# authoritative diff for xl/worksheets/sheet1.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D18").Value = 1324

$ws.Range("D33").Value = 891
$ws.Range("E33").Value = 435

$ws.Range("D34").Value = 940
$ws.Range("E34").Value = 456

$ws.Range("D47").Value = 1292

$ws.Range("D48").Value = 1324
$ws.Range("E48").Value = 390

$ws.Range("D49").Value = 1324
$ws.Range("E49").Value = 345

$ws.Range("D50").Value = 1324

$ws.Range("D51").Value = 507

$ws.Range("E65").Value = 136

$ws.Range("E66").Value = 88

$ws.Range("D67").Value = 127
$ws.Range("E67").Value = 25

$ws.Range("D68").Value = 203
$ws.Range("E68").Value = 25

$ws.Range("D69").Value = 203
$ws.Range("E69").Value = 88

$ws.Range("D70").Value = 203
$ws.Range("E70").Value = 136

$ws.Range("E74").Value = 456

$ws.Range("E91").Value = 88
